$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E)
# figures with the latest scraped values, preserving them as plain text
# (matching the source sheet's inline-string storage of scraped data).
$updates = @{
    "D2" = "261.21"
    "E2" = "1.73%"
    "D3" = "27.41"
    "E3" = "1.21%"
    "E4" = "8.69%"
    "D5" = "0.06070"
    "E5" = "3.07%"
    "D6" = "6.660"
    "D7" = "0.8621"
    "E7" = "0.86%"
    "D8" = "0.9223"
    "E8" = "-1.57%"
    "E9" = "1.40%"
    "D10" = "0.04932"
    "E10" = "6.07%"
    "D11" = "0.07104"
    "E11" = "0.42%"
    "D12" = "0.03049"
    "E12" = "-0.78%"
    "D13" = "0.09085"
    "E13" = "-0.32%"
    "E14" = "-0.63%"
    "D15" = "0.0006108"
    "E15" = "1.29%"
    "D16" = "0.006098"
    "E16" = "-1.00%"
    "D17" = "3.457"
    "E17" = "-1.02%"
    "E18" = "-0.65%"
    "D19" = "2.176"
    "E19" = "-2.25%"
    "E20" = "2.53%"
    "E21" = "2.23%"
    "D22" = "4.112"
    "E22" = "4.96%"
    "D23" = "0.04249"
    "E23" = "-0.71%"
    "E24" = "-0.14%"
    "E25" = "-8.69%"
    "D26" = "0.0001201"
    "E26" = "0.01%"
    "E27" = "3.14%"
    "D40" = "0.03880"
    "E40" = "1.56%"
    "E41" = "1.11%"
    "D42" = "0.004136"
    "E42" = "6.34%"
    "D43" = "0.01502"
    "E43" = "6.63%"
    "E44" = "-7.45%"
    "D45" = "0.00005302"
    "E45" = "-1.07%"
    "E47" = "-17.18%"
    "D48" = "0.1321"
    "E48" = "-47.51%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
